$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Census")

# New "Cleaned" column header (column J)
$ws.Range("J1").Value = "Cleaned"
$ws.Range("J1").Font.Bold = $true

# Mark rows that have already been cleaned with an "x" in the new column
$ws.Range("J5").Value = "x"
$ws.Range("J5").WrapText = $true

$ws.Range("J18").Value = "x"
$ws.Range("J19").Value = "x"
$ws.Range("J20").Value = "x"
$ws.Range("J21").Value = "x"

$ws.Range("J24").Value = "x"
$ws.Range("J24").WrapText = $true

$ws.Range("J26").Value = "x"
$ws.Range("J27").Value = "x"

# Add a note to the poverty-status rows (boundaries changed between 2012-2017)
$note = "Variables aren't consistent between 2012-2017. Don't clean for now, and think about whether we need this table, or use a different table for poverty instead."
$ws.Range("G12").Value = $note
$ws.Range("G12").WrapText = $true
$ws.Range("G13").Value = $note
$ws.Range("G13").WrapText = $true

# Rows need to grow to fit the added note text
$ws.Rows.Item(12).RowHeight = 45
$ws.Rows.Item(13).RowHeight = 45

# Reflect the new scroll position / selection left after the edits
$ws.Activate()
$ws.Range("J23").Select()
